# Auto-generated: apply 2022-06-14 data update to violent-crime-full-year workbook
$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 3008
$ws.Range("I3").Value = 3079
$ws.Range("G4").Value = 1432
$ws.Range("I4").Value = 734
$ws.Range("I6").Value = 3521
$ws.Range("G7").Value = 24654
$ws.Range("I7").Value = 10618

# Sheet 10: Uptown
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I2").Value = 31
$ws.Range("I7").Value = 118

# Sheet 12: Bridgeport
$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("I6").Value = 25
$ws.Range("I7").Value = 56

# Sheet 13: Fuller Park
$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("I6").Value = 12
$ws.Range("I7").Value = 36

# Sheet 15: Woodlawn
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I2").Value = 53
$ws.Range("I3").Value = 65
$ws.Range("I7").Value = 192

# Sheet 16: North Lawndale
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value = 100
$ws.Range("I3").Value = 148
$ws.Range("I6").Value = 146
$ws.Range("I7").Value = 424

# Sheet 19: New City
$ws = $wb.Worksheets.Item('New City')
$ws.Range("I6").Value = 76
$ws.Range("I7").Value = 235

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 99
$ws.Range("I6").Value = 75
$ws.Range("I7").Value = 357
$ws.Range("I8").Value = 669
$ws.Range("I14").Value = 56
$ws.Range("I15").Value = 133
$ws.Range("I18").Value = 74
$ws.Range("I19").Value = 291
$ws.Range("I20").Value = 272
$ws.Range("I26").Value = 12
$ws.Range("I29").Value = 701
$ws.Range("I30").Value = 36
$ws.Range("I33").Value = 491
$ws.Range("I42").Value = 369
$ws.Range("I48").Value = 121
$ws.Range("I50").Value = 48
$ws.Range("I51").Value = 96
$ws.Range("I52").Value = 226
$ws.Range("I53").Value = 120
$ws.Range("I54").Value = 239
$ws.Range("I57").Value = 37
$ws.Range("G63").Value = 197
$ws.Range("I63").Value = 44
$ws.Range("I65").Value = 235
$ws.Range("I67").Value = 424
$ws.Range("I68").Value = 34
$ws.Range("I70").Value = 21
$ws.Range("I71").Value = 30
$ws.Range("I73").Value = 89
$ws.Range("I78").Value = 145
$ws.Range("I82").Value = 11
$ws.Range("I83").Value = 208
$ws.Range("I85").Value = 483
$ws.Range("I87").Value = 16
$ws.Range("I88").Value = 94
$ws.Range("I89").Value = 118
$ws.Range("I90").Value = 127
$ws.Range("I91").Value = 129
$ws.Range("I92").Value = 34
$ws.Range("I93").Value = 61
$ws.Range("I94").Value = 92
$ws.Range("I95").Value = 168
$ws.Range("I99").Value = 192
$ws.Range("G101").Value = 24654
$ws.Range("I101").Value = 10618

# Sheet 20: South Chicago
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I3").Value = 80
$ws.Range("I7").Value = 208

# Sheet 21: West Pullman
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I2").Value = 57
$ws.Range("I7").Value = 168

# Sheet 22: Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 116
$ws.Range("I6").Value = 160
$ws.Range("I7").Value = 491

# Sheet 24: Loop
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I6").Value = 118
$ws.Range("I7").Value = 239

# Sheet 25: Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 217
$ws.Range("I3").Value = 246
$ws.Range("I4").Value = 25
$ws.Range("I7").Value = 701

# Sheet 26: Chatham
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 116
$ws.Range("I3").Value = 80
$ws.Range("I7").Value = 291

# Sheet 28: Lake View
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I2").Value = 17
$ws.Range("I3").Value = 24
$ws.Range("I7").Value = 121

# Sheet 3: South Shore
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 126
$ws.Range("I3").Value = 194
$ws.Range("I4").Value = 25
$ws.Range("I7").Value = 483

# Sheet 30: Ashburn
$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I2").Value = 32
$ws.Range("I7").Value = 75

# Sheet 31: Hermosa
$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("I2").Value = 18
$ws.Range("I3").Value = 16

# Sheet 32: Humboldt Park
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I2").Value = 101
$ws.Range("I3").Value = 124
$ws.Range("I7").Value = 369

# Sheet 35: Rogers Park
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I4").Value = 22
$ws.Range("I6").Value = 57
$ws.Range("I7").Value = 145

# Sheet 40: Washington Park
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I2").Value = 45
$ws.Range("I3").Value = 41
$ws.Range("I7").Value = 129

# Sheet 44: Chicago Lawn
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I3").Value = 79
$ws.Range("I7").Value = 272

# Sheet 45: Calumet Heights
$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("I2").Value = 20
$ws.Range("I7").Value = 74

# Sheet 48: West Lawn
$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("I3").Value = 18
$ws.Range("I7").Value = 61

# Sheet 5: Little Village
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I6").Value = 52
$ws.Range("I7").Value = 226

# Sheet 51: West Loop
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I3").Value = 15
$ws.Range("I7").Value = 92

# Sheet 54: Brighton Park
$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I3").Value = 31
$ws.Range("I7").Value = 133

# Sheet 56: Lincoln Square
$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("I6").Value = 13
$ws.Range("I7").Value = 48

# Sheet 57: East Village
$ws = $wb.Worksheets.Item('East Village')
$ws.Range("I6").Value = 7
$ws.Range("I7").Value = 12

# Sheet 62: Portage Park
$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("I6").Value = 26
$ws.Range("I7").Value = 89

# Sheet 64: Albany Park
$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I2").Value = 33
$ws.Range("I7").Value = 99

# Sheet 66: West Elsdon
$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("I3").Value = 6
$ws.Range("I7").Value = 34

# Sheet 67: O'Hare
$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("I2").Value = 6
$ws.Range("I7").Value = 21

# Sheet 68: United Center
$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I2").Value = 21
$ws.Range("I4").Value = 5
$ws.Range("I7").Value = 94

# Sheet 7: Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 208
$ws.Range("I6").Value = 214
$ws.Range("I7").Value = 669

# Sheet 74: Washington Heights
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I2").Value = 38
$ws.Range("I7").Value = 127

# Sheet 75: Little Italy, UIC
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I6").Value = 44
$ws.Range("I7").Value = 96

# Sheet 76: North Park
$ws = $wb.Worksheets.Item('North Park')
$ws.Range("I3").Value = 11
$ws.Range("I7").Value = 34

# Sheet 77: Mckinley Park
$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("I2").Value = 15
$ws.Range("I7").Value = 37

# Sheet 8: Logan Square
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I6").Value = 52
$ws.Range("I7").Value = 120

# Sheet 81: Oakland
$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("I2").Value = 9
$ws.Range("I7").Value = 30

# Sheet 83: Sheffield & DePaul
$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("I5").Value = 6
$ws.Range("I6").Value = 11

# Sheet 9: Auburn Gresham
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I2").Value = 121
$ws.Range("I6").Value = 94
$ws.Range("I7").Value = 357

# Sheet 92: Ukrainian Village
$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("I3").Value = 5
$ws.Range("I7").Value = 16
